# HTTPS.xlsx GSC export refresh
#
# The "Chart" sheet holds a rolling window of daily rows (Date, Non-HTTPS
# URLs, HTTPS URLs). This update rolls the window forward by one day:
#   - the oldest date (2025-11-06) is dropped
#   - every remaining row shifts up to the next day's figures
#   - a new row for 2026-02-04 is appended at the end
#
# Column A stores the date as plain text (matching the original export),
# so the range is pre-formatted as Text before the values are written --
# otherwise Excel would helpfully "autocorrect" a typed date-shaped string
# like "2025-11-07" into a real date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$firstRow = 2
$lastRow = 91

$ws.Range("A" + $firstRow + ":A" + $lastRow).NumberFormat = "@"

# New Date column values (row 2 .. row 91), already shifted forward one day.
$dates = @(
    "2025-11-07",
    "2025-11-08",
    "2025-11-09",
    "2025-11-10",
    "2025-11-11",
    "2025-11-12",
    "2025-11-13",
    "2025-11-14",
    "2025-11-15",
    "2025-11-16",
    "2025-11-17",
    "2025-11-18",
    "2025-11-19",
    "2025-11-20",
    "2025-11-21",
    "2025-11-22",
    "2025-11-23",
    "2025-11-24",
    "2025-11-25",
    "2025-11-26",
    "2025-11-27",
    "2025-11-28",
    "2025-11-29",
    "2025-11-30",
    "2025-12-01",
    "2025-12-02",
    "2025-12-03",
    "2025-12-04",
    "2025-12-05",
    "2025-12-06",
    "2025-12-07",
    "2025-12-08",
    "2025-12-09",
    "2025-12-10",
    "2025-12-11",
    "2025-12-12",
    "2025-12-13",
    "2025-12-14",
    "2025-12-15",
    "2025-12-16",
    "2025-12-17",
    "2025-12-18",
    "2025-12-19",
    "2025-12-20",
    "2025-12-21",
    "2025-12-22",
    "2025-12-23",
    "2025-12-24",
    "2025-12-25",
    "2025-12-26",
    "2025-12-27",
    "2025-12-28",
    "2025-12-29",
    "2025-12-30",
    "2025-12-31",
    "2026-01-01",
    "2026-01-02",
    "2026-01-03",
    "2026-01-04",
    "2026-01-05",
    "2026-01-06",
    "2026-01-07",
    "2026-01-08",
    "2026-01-09",
    "2026-01-10",
    "2026-01-11",
    "2026-01-12",
    "2026-01-13",
    "2026-01-14",
    "2026-01-15",
    "2026-01-16",
    "2026-01-17",
    "2026-01-18",
    "2026-01-19",
    "2026-01-20",
    "2026-01-21",
    "2026-01-22",
    "2026-01-23",
    "2026-01-24",
    "2026-01-25",
    "2026-01-26",
    "2026-01-27",
    "2026-01-28",
    "2026-01-29",
    "2026-01-30",
    "2026-01-31",
    "2026-02-01",
    "2026-02-02",
    "2026-02-03",
    "2026-02-04"
)

# New "HTTPS URLs" (column C) values for the same rows.
$pages = @(
    94.0,
    86.0,
    83.0,
    66.0,
    54.0,
    46.0,
    43.0,
    40.0,
    37.0,
    35.0,
    30.0,
    29.0,
    26.0,
    25.0,
    25.0,
    26.0,
    26.0,
    25.0,
    25.0,
    27.0,
    28.0,
    28.0,
    27.0,
    27.0,
    27.0,
    27.0,
    27.0,
    26.0,
    25.0,
    25.0,
    25.0,
    26.0,
    27.0,
    27.0,
    29.0,
    29.0,
    30.0,
    30.0,
    31.0,
    31.0,
    31.0,
    31.0,
    31.0,
    32.0,
    32.0,
    32.0,
    32.0,
    30.0,
    31.0,
    32.0,
    30.0,
    28.0,
    28.0,
    28.0,
    28.0,
    29.0,
    29.0,
    28.0,
    27.0,
    27.0,
    28.0,
    27.0,
    27.0,
    27.0,
    27.0,
    26.0,
    26.0,
    27.0,
    26.0,
    26.0,
    25.0,
    25.0,
    25.0,
    25.0,
    26.0,
    25.0,
    24.0,
    23.0,
    24.0,
    24.0,
    24.0,
    25.0,
    26.0,
    27.0,
    28.0,
    28.0,
    28.0,
    28.0,
    28.0,
    28.0
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 3).Value = $pages[$i]
}
